$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Row 24: was a lone "did I put new selection on too?" note in D24 only.
# It becomes a full log row: Neurolucida results / DEG run with the new BL_A+BL_C selection.
$ws.Range("A24").Value = "Neurolucida results"
$ws.Range("B24").Value = "2022-06-13 16-27-52"
$ws.Range("C24").Value = "DEG"
$ws.Range("D24").Value = "SCTv2 corrected BL_A + BL_C new selection"

# Row 25 and Row 26 stay the same content-wise (only the shared-string table
# shrinks because the old D24 note is gone, which the runtime reindexes for us).

# New row 27: Neurolucida / DEG run for BL_N + BL_C old post selection.
$ws.Range("A27").Value = "Neurolucida results"
$ws.Range("C27").Value = "DEG"
$ws.Range("D27").Value = "SCTv2 corrected BL_N + BL_C old post selection"
$ws.Range("B27").Value = "2022-06-14 14-45-40"
$ws.Range("F27").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G27").Value = "pseudotime"

# New row 28: Neurolucida / DEG run for BL_N + BL_C old selection.
$ws.Range("A28").Value = "Neurolucida results"
$ws.Range("B28").Value = "2022-06-14 14-46-19"
$ws.Range("C28").Value = "DEG"
$ws.Range("D28").Value = "SCTv2 corrected BL_N + BL_C old selection"
$ws.Range("F28").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G28").Value = "pseudotime"

# Update the saved selection to match where the author left off editing.
$ws.Range("F32").Select()
